# Generate Report for Handoff
# Adds a new tracked file "bc575ac8-ea61-43c4-a217-9a7e9f6bdfa8" as row 9
# to the "Overview", "zh-cn" and "de-de" worksheets of the localization
# status report, mirroring the existing "Ready for handoff" rows.

$wb = $excel.ActiveWorkbook

$fileBase = "bc575ac8-ea61-43c4-a217-9a7e9f6bdfa8"
$mdName   = "$fileBase.md"
$xlfHash  = "4c23f776e6311b8bc14215dd8d98839de6760442"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"  (columns: File Name | zh-cn | de-de | Latest Handoff Date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsOverview.Range("B9").Value = "Ready for handoff"
$wsOverview.Range("C9").Value = "Ready for handoff"
$wsOverview.Range("D9").Value = "2016-03-25 07:27:21"
$wsOverview.Range("D9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn" (columns: Source File Name | File Extension | Status |
#   Latest Handoff File | Latest Handoff Datetime | Latest Target File |
#   Latest Handback File | Latest Handback DateTime | Reference Tokens |
#   Handoff Reason | Dependency From | Error Detail)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhXlf = "$fileBase.$xlfHash.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf",
    "",
    "",
    $zhXlf
)
$wsZhCn.Range("E9").Value = "2016-03-25 07:27:16"
$wsZhCn.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("J9").Value = "Include"

# ---------------------------------------------------------------------
# Sheet 3: "de-de" (same column layout as "zh-cn")
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$deXlf = "$fileBase.$xlfHash.de-de.xlf"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf",
    "",
    "",
    $deXlf
)
$wsDeDe.Range("E9").Value = "2016-03-25 07:27:21"
$wsDeDe.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("J9").Value = "Include"

"Handoff report row added for $fileBase"
